# Insert a new weekly price record as row 636 ("Fruta / hortaliza, semanal").
# This pushes the existing rows 636-683 down to 637-684 (dimension grows to
# A1:T684) and populates the newly inserted row with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(636).Insert()

$ws.Cells.Item(636, 1).Value  = 6
$ws.Cells.Item(636, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(636, 3).Value  = "Metropolitana"
$ws.Cells.Item(636, 4).Value  = 45132
$ws.Cells.Item(636, 5).Value  = 13
$ws.Cells.Item(636, 6).Value  = "Fruta"
$ws.Cells.Item(636, 7).Value  = 100101
$ws.Cells.Item(636, 8).Value  = "Berries"
$ws.Cells.Item(636, 9).Value  = 100101001
$ws.Cells.Item(636, 10).Value = "Arándano (blue)"
$ws.Cells.Item(636, 11).Value = "Sin especificar"
$ws.Cells.Item(636, 12).Value = "Primera"
$ws.Cells.Item(636, 13).Value = 730
$ws.Cells.Item(636, 14).Value = 14000
$ws.Cells.Item(636, 15).Value = 14000
$ws.Cells.Item(636, 16).Value = 14000
$ws.Cells.Item(636, 17).Value = "$/bandeja 12 canastillos 125 gramos"
$ws.Cells.Item(636, 18).Value = "Perú"
$ws.Cells.Item(636, 19).Value = 9333
$ws.Cells.Item(636, 20).Value = 1.5
